$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the manifest/data file name referenced in C2 (was "...Neo4jData.xlsx", now "...Manifest.csv")
$ws.Range("C2").Value = "TC01_Bento_E2E_Select-All-Add-To-Cart_Manifest.csv"

# Update the view: drop the top-left-cell scroll position and set zoom to 70%
$win = $excel.ActiveWindow
$win.Zoom = 70
$win.ScrollRow = 1
$win.ScrollColumn = 1

$ws.Range("B2").Select()
